$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the duty-cycle parameters for rows 2 and 3 (30 on / 30 off continuous pulse)
$ws.Range("C2").Value = 30000
$ws.Range("D2").Value = 30000
$ws.Range("E2").Value = 1

$ws.Range("C3").Value = 30000
$ws.Range("D3").Value = 30000
$ws.Range("E3").Value = 1

# Update the selected cell in the sheet view
$ws.Activate()
$ws.Range("D9").Select()
